$wb = $excel.ActiveWorkbook

# --- Metadata sheet ---
$meta = $wb.Worksheets.Item("Metadata")

# The sheet currently has a duplicated "Contact" / "No display for ContactDetail" row
# (rows 10 and 11 are identical). Remove one of them so everything below shifts up
# by one row (new last row becomes 20 instead of 21).
$meta.Rows.Item(11).Delete()

# Version bump
$meta.Range("B3").Value = "6.0.0"

# Date bump
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value
$meta.Range("B9").Value = "Alvearie Team"

# The now-single former "Contact" row (row 10) becomes the new "Jurisdiction" row
$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"

# --- Elements sheet ---
$elements = $wb.Worksheets.Item("Elements")

# Root Extension element's Short/Definition now reflect the profile name/description
# instead of the generic Extension placeholders.
$elements.Range("K2").Value = "Union Worker Indicator"
$elements.Range("L2").Value = "Indicator that the contract holder (subscriber) belongs to a union"
